$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "actual hours burned" tally for rows 32-37 from column M to column N
for ($r = 32; $r -le 37; $r++) {
    $ws.Cells.Item($r, 13).Clear()              # M32:M37 -> clear (cell removed entirely)
    $ws.Cells.Item($r, 14).Value = 1            # N32:N37 -> 1
}

# F15 used to be part of the shared formula "D15-SUM(L25:L37)" style chain;
# give it its own explicit formula pulling from the now-used column M range.
$ws.Range("F15").Formula = "=E15-SUM(M25:M37)"

# G15 likewise becomes explicit, now referencing the column N range instead
# of the old shared-formula pattern.
$ws.Range("G15").Formula = "=F15-SUM(N25:N37)"

$wb.Save()
